$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $srcRow = 45
    $dstRow = 46

    # Column A: date/time value, copy the number format from the row above.
    $ws.Cells.Item($dstRow, 1).Value2 = 45832.46327546296
    $ws.Cells.Item($dstRow, 1).NumberFormat = $ws.Cells.Item($srcRow, 1).NumberFormat

    # Columns B through I: duplicate the values straight across.
    for ($col = 2; $col -le 9; $col++) {
        $srcCell = $ws.Cells.Item($srcRow, $col)
        $dstCell = $ws.Cells.Item($dstRow, $col)
        $dstCell.Value2 = $srcCell.Value2
    }
}
